# Refresh the crypto price/volume snapshot (cols D & E) for rows 2-51,
# plus a two-row reorder at 44/45 (FraxShare <-> HuobiToken) per the upstream feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-looking numbers (e.g. '228.49') but the source sheet
# stores them as TEXT (inline strings), not numbers - price formatting (and
# trailing zeros like '22.00') must be preserved verbatim. A leading apostrophe
# forces Excel to keep the literal text instead of coercing it to a Double.

$ws.Range("D2").Value = '38.705.32'
$ws.Range("E2").Value = '  +1.01%  '

$ws.Range("D3").Value = '2.099.27'
$ws.Range("E3").Value = '  +0.10%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = '''228.49'
$ws.Range("E5").Value = '  -0.09%  '

$ws.Range("D6").Value = '''0.616'
$ws.Range("E6").Value = '  +0.69%  '

$ws.Range("D7").Value = '''61.93'
$ws.Range("E7").Value = '  +1.11%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("D9").Value = '''0.387'
$ws.Range("E9").Value = '  +1.79%  '

$ws.Range("D10").Value = '''0.0840'
$ws.Range("E10").Value = '  -0.51%  '

$ws.Range("E11").Value = '  -0.18%  '

$ws.Range("D12").Value = '''15.84'
$ws.Range("E12").Value = '  +7.25%  '

$ws.Range("D13").Value = '2.410.37'
$ws.Range("E13").Value = '  +0.04%  '

$ws.Range("D14").Value = '''22.00'
$ws.Range("E14").Value = '  -1.62%  '

$ws.Range("D15").Value = '''0.804'
$ws.Range("E15").Value = '  +3.49%  '

$ws.Range("D16").Value = '''5.53'
$ws.Range("E16").Value = '  +1.42%  '

$ws.Range("D17").Value = '2.086.85'
$ws.Range("E17").Value = '  -0.80%  '

$ws.Range("D18").Value = '38.729.83'
$ws.Range("E18").Value = '  +1.11%  '

$ws.Range("D19").Value = '''71.84'
$ws.Range("E19").Value = '  +1.97%  '

$ws.Range("D20").Value = '''6.06'
$ws.Range("E20").Value = '  +0.95%  '

$ws.Range("D21").Value = '0.0₃0839'
$ws.Range("E21").Value = '  +0.69%  '

$ws.Range("D22").Value = '''227.63'
$ws.Range("E22").Value = '  +1.04%  '

$ws.Range("E24").Value = '  -1.83%  '

$ws.Range("D25").Value = '''2.33'
$ws.Range("E25").Value = '  +0.27%  '

$ws.Range("D26").Value = '''172.22'
$ws.Range("E26").Value = '  +1.22%  '

$ws.Range("D27").Value = '''9.55'
$ws.Range("E27").Value = '  +1.60%  '

$ws.Range("E28").Value = '  +6.02%  '

$ws.Range("E29").Value = '  +4.56%  '

$ws.Range("E30").Value = '  +1.38%  '

$ws.Range("D31").Value = '''2.51'
$ws.Range("E31").Value = '  +6.95%  '

$ws.Range("E32").Value = '  +0.67%  '

$ws.Range("E33").Value = '  +2.48%  '

$ws.Range("E34").Value = '  +0.27%  '

$ws.Range("D35").Value = '''0.0619'
$ws.Range("E35").Value = '  +2.54%  '

$ws.Range("E36").Value = '  +5.57%  '

$ws.Range("D37").Value = '''2.42'
$ws.Range("E37").Value = '  +1.27%  '

$ws.Range("D38").Value = '''3.57'
$ws.Range("E38").Value = '  +1.60%  '

$ws.Range("E39").Value = '  -0.16%  '

$ws.Range("D40").Value = '''18.16'
$ws.Range("E40").Value = '  -0.33%  '

$ws.Range("D41").Value = '''0.0228'
$ws.Range("E41").Value = '  +3.96%  '

$ws.Range("E42").Value = '  +2.49%  '

$ws.Range("D43").Value = '1.533.64'
$ws.Range("E43").Value = '  -0.76%  '

$ws.Range("B44").Value = 'HuobiToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D44").Value = '''2.80'
$ws.Range("E44").Value = '  -1.16%  '

$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = '''7.83'
$ws.Range("E45").Value = '  +4.18%  '

$ws.Range("E46").Value = '  +2.79%  '

$ws.Range("D47").Value = '''0.0910'
$ws.Range("E47").Value = '  -0.05%  '

$ws.Range("E48").Value = '  -0.11%  '

$ws.Range("D49").Value = '''1.04'
$ws.Range("E49").Value = '  +1.30%  '

$ws.Range("E50").Value = '  -0.47%  '

$ws.Range("D51").Value = '2.292.12'
$ws.Range("E51").Value = '  -0.09%  '
